$d = $word.ActiveDocument

# Region A: "<<{dateFormat(creationDate" -- merges split '<' + '<{' runs into one run
# and drops the stray gramStart/gramEnd proofErr markers picked up along the way.
$d.Content.Find.Execute("<<{dateFormat(creationDate", $false, $false, $false, $false, $false, `
    $true, 1, $false, "<<{dateFormat(creationDate", 2) | Out-Null

# Region B: "<<{dateFormat(hearingDate" -- only the gramStart/gramEnd markers are dropped,
# the literal text is unchanged.
$d.Content.Find.Execute("<<{dateFormat(hearingDate", $false, $false, $false, $false, $false, `
    $true, 1, $false, "<<{dateFormat(hearingDate", 2) | Out-Null

# Region D: drop the leading "by " before the <<hearingType>> field.
$d.Content.Find.Execute("by <<hearingType", $false, $false, $false, $false, $false, `
    $true, 1, $false, "<<hearingType", 2) | Out-Null

# Region C: "<<cs_{additionalInfo!=null}>>" -- merges the "!=" and "null}>>" runs and
# removes the gramStart/gramEnd markers in between.
$d.Content.Find.Execute("<<cs_{additionalInfo!=null}>>", $false, $false, $false, $false, $false, `
    $true, 1, $false, "<<cs_{additionalInfo!=null}>>", 2) | Out-Null

# Region E: "Payable by <<{dateFormat(hearingDueDate" -- merges the '<' + '<{' runs and
# drops the gramStart/gramEnd markers.
$d.Content.Find.Execute("Payable by <<{dateFormat(hearingDueDate", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Payable by <<{dateFormat(hearingDueDate", 2) | Out-Null
